$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in the existing "html and DRY principle??" checklist item.
$ws.Range("A25").Value = "html and DRY principle?? (in its defense i'll say that html is not a programming language after all)"

# Mark "when to call dto and when to call model" (A24) as done, same green
# fill used by the other completed checklist items.
$ws.Range("A24").Interior.Color = 5287936

# Append the three new checklist items.
$ws.Range("A26").Value = "how to edit log in page"
$ws.Range("A27").Value = "don't like async everywhere in service layer"
$ws.Range("A28").Value = "[.. ]"

# Move the selection to the new last row, matching where the author left off.
$ws.Range("A27").Select()
